# Update "想去人数" (number of people interested) figures that changed
# between data refreshes for the 南宁·2024三月三国潮动漫节（良牙春典）
# and 南宁·布谷鸟动漫展4th events.
#
# These two rows appear in both the "展览" sheet and the "全部类型" sheet,
# so both need to be updated.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 5454
    $ws.Range("F4").Value = 936
}
